$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.972.70'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '4.030.61'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  -0.08%  '
$origStyle1 = $ws.Range('D5').Style
$ws.Range('D5').Value = "'535.73"
$ws.Range('D5').Style = $origStyle1
$ws.Range('E5').Value = '  +1.16%  '
$origStyle2 = $ws.Range('D6').Style
$ws.Range('D6').Value = "'149.63"
$ws.Range('D6').Style = $origStyle2
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('D7').Value = '4.023.58'
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -2.38%  '
$origStyle3 = $ws.Range('D11').Style
$ws.Range('D11').Value = "'0.173"
$ws.Range('D11').Style = $origStyle3
$ws.Range('E11').Value = '  -2.58%  '
$origStyle4 = $ws.Range('D12').Style
$ws.Range('D12').Value = "'54.01"
$ws.Range('D12').Style = $origStyle4
$ws.Range('E12').Value = '  +6.64%  '
$ws.Range('E13').Value = '  -2.65%  '
$origStyle5 = $ws.Range('D14').Style
$ws.Range('D14').Value = "'10.82"
$ws.Range('D14').Style = $origStyle5
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').Value = '4.667.79'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '4.036.46'
$ws.Range('E16').Value = '  -0.41%  '
$origStyle6 = $ws.Range('D17').Style
$ws.Range('D17').Value = "'14.20"
$ws.Range('D17').Style = $origStyle6
$ws.Range('E17').Value = '  -1.89%  '
$origStyle7 = $ws.Range('D18').Style
$ws.Range('D18').Value = "'20.89"
$ws.Range('D18').Style = $origStyle7
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('E19').Value = '  -3.32%  '
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '71.922.22'
$ws.Range('E21').Value = '  -0.35%  '
$origStyle8 = $ws.Range('D22').Style
$ws.Range('D22').Value = "'431.56"
$ws.Range('D22').Style = $origStyle8
$ws.Range('E22').Value = '  -1.72%  '
$origStyle9 = $ws.Range('D23').Style
$ws.Range('D23').Value = "'98.31"
$ws.Range('D23').Style = $origStyle9
$ws.Range('E23').Value = '  -2.50%  '
$origStyle10 = $ws.Range('D24').Style
$ws.Range('D24').Value = "'3.61"
$ws.Range('D24').Style = $origStyle10
$ws.Range('E24').Value = '  -2.13%  '
$origStyle11 = $ws.Range('D25').Style
$ws.Range('D25').Value = "'14.79"
$ws.Range('D25').Style = $origStyle11
$ws.Range('E25').Value = '  -2.08%  '
$origStyle12 = $ws.Range('D26').Style
$ws.Range('D26').Value = "'4.23"
$ws.Range('D26').Style = $origStyle12
$ws.Range('E26').Value = '  -0.45%  '
$origStyle13 = $ws.Range('D27').Style
$ws.Range('D27').Value = "'4.36"
$ws.Range('D27').Style = $origStyle13
$ws.Range('E27').Value = '  +28.29%  '
$origStyle14 = $ws.Range('D28').Style
$ws.Range('D28').Value = "'11.41"
$ws.Range('D28').Style = $origStyle14
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('E29').Value = '  -2.32%  '
$origStyle15 = $ws.Range('D30').Style
$ws.Range('D30').Value = "'5.95"
$ws.Range('D30').Style = $origStyle15
$ws.Range('E30').Value = '  +1.80%  '
$origStyle16 = $ws.Range('D31').Style
$ws.Range('D31').Value = "'37.10"
$ws.Range('D31').Style = $origStyle16
$ws.Range('E31').Value = '  -1.20%  '
$origStyle17 = $ws.Range('D32').Style
$ws.Range('D32').Value = "'8.32"
$ws.Range('D32').Style = $origStyle17
$ws.Range('E32').Value = '  +22.30%  '
$ws.Range('E33').Value = '  +2.17%  '
$origStyle18 = $ws.Range('D34').Style
$ws.Range('D34').Value = "'50.28"
$ws.Range('D34').Style = $origStyle18
$ws.Range('E34').Value = '  +17.28%  '
$origStyle19 = $ws.Range('D35').Style
$ws.Range('D35').Value = "'13.62"
$ws.Range('D35').Style = $origStyle19
$ws.Range('E35').Value = '  -1.12%  '
$origStyle20 = $ws.Range('D36').Style
$ws.Range('D36').Value = "'676.04"
$ws.Range('D36').Style = $origStyle20
$ws.Range('E36').Value = '  -0.06%  '
$origStyle21 = $ws.Range('D37').Style
$ws.Range('D37').Value = "'68.14"
$ws.Range('D37').Style = $origStyle21
$ws.Range('E37').Value = '  +2.01%  '
$origStyle22 = $ws.Range('D38').Style
$ws.Range('D38').Value = "'0.461"
$ws.Range('D38').Style = $origStyle22
$ws.Range('E38').Value = '  +4.26%  '
$ws.Range('E39').Value = '  -5.62%  '
$ws.Range('E40').Value = '  +8.26%  '
$ws.Range('E41').Value = '  -6.00%  '
$origStyle23 = $ws.Range('D42').Style
$ws.Range('D42').Value = "'3.41"
$ws.Range('D42').Style = $origStyle23
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('B43').Value = 'Dai'
$ws.Range('C43').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$origStyle24 = $ws.Range('D43').Style
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = $origStyle24
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('B44').Value = 'THORChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$origStyle25 = $ws.Range('D44').Style
$ws.Range('D44').Value = "'11.05"
$ws.Range('D44').Style = $origStyle25
$ws.Range('E44').Value = '  +15.23%  '
$origStyle26 = $ws.Range('D45').Style
$ws.Range('D45').Value = "'0.0495"
$ws.Range('D45').Style = $origStyle26
$ws.Range('E45').Value = '  -2.41%  '
$origStyle27 = $ws.Range('D46').Style
$ws.Range('D46').Value = "'1.00"
$ws.Range('D46').Style = $origStyle27
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  -2.37%  '
$origStyle28 = $ws.Range('D48').Style
$ws.Range('D48').Value = "'2.67"
$ws.Range('D48').Style = $origStyle28
$ws.Range('E48').Value = '  -4.79%  '
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').Value = '2.875.76'
$ws.Range('E51').Value = '  +9.92%  '
